$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the RF (column I) values for rows 31-54 to reflect the 2025 data / RF change
$ws.Range("I31:I54").Value = 17.87725806451613
